$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 6711.1
$ws.Range("J40").Value = 7206.5713
$ws.Range("L40").Value = 7206.5713
$ws.Range("N40").Value = -7556.5713
$ws.Range("H51").Value = 0
$ws.Range("J51").Value = 0
$ws.Range("L51").Value = 0
$ws.Range("N51").Value = ""
$ws.Range("H70").Value = 3405.8
$ws.Range("I70").Value = 982.8333
$ws.Range("J70").Value = 5021.1113
$ws.Range("K70").Value = 2948.4999
$ws.Range("L70").Value = 15063.3339
$ws.Range("M70").Value = -2678.4999
$ws.Range("N70").Value = -15603.3339
$ws.Range("H73").Value = 3405.8
$ws.Range("I73").Value = 982.8333
$ws.Range("J73").Value = 5021.1113
$ws.Range("K73").Value = 2948.4999
$ws.Range("L73").Value = 15063.3339
$ws.Range("M73").Value = -2012.4999
$ws.Range("N73").Value = -16935.3339
$ws.Range("H76").Value = 4166.5
$ws.Range("I76").Value = 3333
$ws.Range("K76").Value = 3333
$ws.Range("M76").Value = -3018
$ws.Range("H79").Value = 4166.5
$ws.Range("I79").Value = 3333
$ws.Range("K79").Value = 3333
$ws.Range("M79").Value = -2241
$ws.Range("H98").Value = 731.1667
$ws.Range("I98").Value = 656.8
$ws.Range("J98").Value = 1103
$ws.Range("K98").Value = 656.8
$ws.Range("L98").Value = 1103
$ws.Range("M98").Value = 841.2
$ws.Range("N98").Value = -4099
$ws.Range("H103").Value = 4737.1763
$ws.Range("J103").Value = 5235.6924
$ws.Range("L103").Value = 15707.0772
$ws.Range("N103").Value = -16879.0772
$ws.Range("H116").Value = 10709.333
$ws.Range("I116").Value = 10798
$ws.Range("K116").Value = 10798
$ws.Range("M116").Value = -7356
$ws.Range("H122").Value = 731.1667
$ws.Range("I122").Value = 656.8
$ws.Range("J122").Value = 1103
$ws.Range("K122").Value = 1970.4
$ws.Range("L122").Value = 3309
$ws.Range("M122").Value = 479.6000000000001
$ws.Range("N122").Value = -8209

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 3554.9443
$ws.Range("I61").Value = 2332.5
$ws.Range("K61").Value = 2332.5
$ws.Range("M61").Value = -2120.5
$ws.Range("H74").Value = 6419
$ws.Range("I74").Value = 6000
$ws.Range("K74").Value = 6000
$ws.Range("M74").Value = -5126
$ws.Range("H77").Value = 6419
$ws.Range("I77").Value = 6000
$ws.Range("K77").Value = 30000
$ws.Range("M77").Value = -25632
$ws.Range("H136").Value = 3554.9443
$ws.Range("I136").Value = 2332.5
$ws.Range("K136").Value = 6997.5
$ws.Range("M136").Value = -4447.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 5224.4287
$ws.Range("I86").Value = 1642.75
$ws.Range("K86").Value = 1642.75
$ws.Range("M86").Value = -519.75
$ws.Range("H89").Value = 5224.4287
$ws.Range("I89").Value = 1642.75
$ws.Range("K89").Value = 8213.75
$ws.Range("M89").Value = -2597.75
$ws.Range("H105").Value = 8265795.5
$ws.Range("I105").Value = 18182990
$ws.Range("K105").Value = 18182990
$ws.Range("M105").Value = -18181243
$ws.Range("H107").Value = 35719084
$ws.Range("I107").Value = 62501836
$ws.Range("K107").Value = 62501836
$ws.Range("M107").Value = -62499916

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 58.857143
$ws.Range("I7").Value = 49.2
$ws.Range("K7").Value = 49.2
$ws.Range("M7").Value = 63.8
$ws.Range("H11").Value = 1250
$ws.Range("I11").Value = 500
$ws.Range("J11").Value = 5000
$ws.Range("K11").Value = 500
$ws.Range("L11").Value = 5000
$ws.Range("M11").Value = -360
$ws.Range("N11").Value = -5280
$ws.Range("H14").Value = 3116.6667
$ws.Range("I14").Value = 2675
$ws.Range("K14").Value = 2675
$ws.Range("M14").Value = -2505
$ws.Range("H17").Value = 3333.3333
$ws.Range("I17").Value = 2000
$ws.Range("K17").Value = 2000
$ws.Range("M17").Value = -1826
$ws.Range("H99").Value = 4175.3335
$ws.Range("J99").Value = 5257
$ws.Range("L99").Value = 5257
$ws.Range("N99").Value = -8253
$ws.Range("H105").Value = 1315.6
$ws.Range("I105").Value = 1239.5555
$ws.Range("K105").Value = 1239.5555
$ws.Range("M105").Value = 507.4445000000001
$ws.Range("H126").Value = 4175.3335
$ws.Range("J126").Value = 5257
$ws.Range("L126").Value = 15771
$ws.Range("N126").Value = -20711
$ws.Range("H132").Value = 2284.4783
$ws.Range("I132").Value = 2113.8333
$ws.Range("K132").Value = 6341.499899999999
$ws.Range("M132").Value = -3811.499899999999
$ws.Range("H134").Value = 2223.739
$ws.Range("I134").Value = 1057.375
$ws.Range("K134").Value = 3172.125
$ws.Range("M134").Value = -637.125

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H11").Value = 16926
$ws.Range("J11").Value = 3124.5
$ws.Range("L11").Value = 9373.5
$ws.Range("N11").Value = -9653.5
$ws.Range("H32").Value = 811.375
$ws.Range("J32").Value = 755.8570999999999
$ws.Range("L32").Value = 2267.5713
$ws.Range("N32").Value = -2833.5713
$ws.Range("H121").Value = 556.4286
$ws.Range("J121").Value = 639.2
$ws.Range("L121").Value = 1917.6
$ws.Range("N121").Value = -4537.6
$ws.Range("H129").Value = 2229
$ws.Range("I129").Value = 1100
$ws.Range("J129").Value = 2605.3333
$ws.Range("K129").Value = 3300
$ws.Range("L129").Value = 7815.999899999999
$ws.Range("M129").Value = 1700
$ws.Range("N129").Value = -17815.9999
$ws.Range("H131").Value = 1644
$ws.Range("I131").Value = 1192.5834
$ws.Range("K131").Value = 3577.7502
$ws.Range("M131").Value = 1462.2498
$ws.Range("H134").Value = 3624.5
$ws.Range("I134").Value = 3624.5
$ws.Range("K134").Value = 10873.5
$ws.Range("M134").Value = -5803.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 287.0625
$ws.Range("I2").Value = 160.77777
$ws.Range("J2").Value = 449.42856
$ws.Range("K2").Value = 160.77777
$ws.Range("L2").Value = 449.42856
$ws.Range("M2").Value = -47.77777
$ws.Range("N2").Value = -675.4285600000001
$ws.Range("H70").Value = 0
$ws.Range("I70").Value = 0
$ws.Range("K70").Value = 0
$ws.Range("M70").Value = ""
$ws.Range("H73").Value = 0
$ws.Range("I73").Value = 0
$ws.Range("K73").Value = 0
$ws.Range("M73").Value = ""
$ws.Range("H113").Value = 9249.4
$ws.Range("I113").Value = 8123.5
$ws.Range("K113").Value = 8123.5
$ws.Range("M113").Value = -5953.5
$ws.Range("H126").Value = 8066.1665
$ws.Range("I126").Value = 6599.3335
$ws.Range("K126").Value = 19798.0005
$ws.Range("M126").Value = -17328.0005
$ws.Range("H140").Value = 85965
$ws.Range("I140").Value = 0
$ws.Range("J140").Value = 85965
$ws.Range("K140").Value = 0
$ws.Range("L140").Value = 85965
$ws.Range("M140").Value = ""
$ws.Range("N140").Value = -96325

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H2").Value = 2933.3333
$ws.Range("H22").Value = 1709.2
$ws.Range("J22").Value = 1741.7142
$ws.Range("L22").Value = 1741.7142
$ws.Range("N22").Value = -2331.7142
$ws.Range("H27").Value = 1709.2
$ws.Range("J27").Value = 1741.7142
$ws.Range("L27").Value = 1741.7142
$ws.Range("N27").Value = -1955.7142
$ws.Range("H43").Value = 0
$ws.Range("I43").Value = 0
$ws.Range("K43").Value = 0
$ws.Range("M43").Value = ""
$ws.Range("H46").Value = 847.8
$ws.Range("I46").Value = 700
$ws.Range("J46").Value = 884.75
$ws.Range("K46").Value = 700
$ws.Range("L46").Value = 884.75
$ws.Range("M46").Value = -512
$ws.Range("N46").Value = -1260.75
$ws.Range("H120").Value = 47666.668
$ws.Range("J120").Value = 47666.668
$ws.Range("L120").Value = 47666.668
$ws.Range("N120").Value = -57342.668
$ws.Range("H122").Value = 916.3333
$ws.Range("I122").Value = 916.3333
$ws.Range("K122").Value = 2748.9999
$ws.Range("M122").Value = -298.9998999999998

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H4").Value = 51071.145
$ws.Range("I4").Value = 58333.332
$ws.Range("J4").Value = 7498
$ws.Range("K4").Value = 58333.332
$ws.Range("L4").Value = 7498
$ws.Range("M4").Value = -58220.332
$ws.Range("N4").Value = -7724
$ws.Range("H51").Value = 23894
$ws.Range("I51").Value = 23894
$ws.Range("K51").Value = 23894
$ws.Range("M51").Value = -23384
$ws.Range("H96").Value = 2079.4167
$ws.Range("I96").Value = 2808.4285
$ws.Range("J96").Value = 1058.8
$ws.Range("K96").Value = 2808.4285
$ws.Range("L96").Value = 1058.8
$ws.Range("M96").Value = -1435.4285
$ws.Range("N96").Value = -3804.8
$ws.Range("H120").Value = 100420
$ws.Range("J120").Value = 100420
$ws.Range("L120").Value = 100420
$ws.Range("N120").Value = -110096
$ws.Range("H132").Value = 2210.7778
$ws.Range("I132").Value = 1725
$ws.Range("J132").Value = 2599.4
$ws.Range("K132").Value = 5175
$ws.Range("L132").Value = 7798.200000000001
$ws.Range("M132").Value = -2645
$ws.Range("N132").Value = -12858.2
